$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Save the existing "species codes" comment text (currently in B11) so it can
# be moved down to the new row 12.
$speciesCodesComment = $ws.Range("B11").Value2

# Row 12 gets the ID 11 and inherits the comment that used to live in row 11.
$ws.Range("A12").Value = 11
$ws.Range("B12").Value = $speciesCodesComment

# Row 11 now holds the new "minimum size" comment (Boutprocess bug fix).
$ws.Range("B11").Value = "If desired, enter minimum size (mm) for items of this prey type "

# Match the saved selection state from the edited workbook.
$ws.Range("B19").Select()
